$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (closest values reachable through the ColumnWidth
# COM property, which Excel itself quantizes to whole pixels).
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334

# Update cell values
$ws.Range("A1").Value = 0.058292534492879806
$ws.Range("B1").Value = -0.058292535669159075
$ws.Range("A2").Value = -0.00055241910279379101
$ws.Range("B2").Value = 0.00055241787173139403
$ws.Range("A3").Value = -0.045310345580090339
$ws.Range("B3").Value = 0.045310344373800292
